$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.190.38'
$ws.Range("E2").Value = '  +0.41%  '

# Row 3
$ws.Range("D3").Value = '2.519.10'
$ws.Range("E3").Value = '  +0.63%  '

# Row 4
$ws.Range("E4").Value = '  +0.07%  '

# Row 5
$ws.Range("D5").Value = '536.44'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6
$ws.Range("D6").Value = '140.48'
$ws.Range("E6").Value = '  -2.59%  '

# Row 7
$ws.Range("E7").Value = '  +0.28%  '

# Row 8
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  -1.08%  '

# Row 9
$ws.Range("D9").Value = '2.523.86'
$ws.Range("E9").Value = '  -0.13%  '

# Row 10
$ws.Range("E10").Value = '  +0.16%  '

# Row 11
$ws.Range("E11").Value = '  +0.91%  '

# Row 12
$ws.Range("D12").Value = '5.48'
$ws.Range("E12").Value = '  -2.74%  '

# Row 13
$ws.Range("D13").Value = '0.359'
$ws.Range("E13").Value = '  +2.20%  '

# Row 14
$ws.Range("D14").Value = '2.966.13'
$ws.Range("E14").Value = '  +0.83%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").Value = '23.07'
$ws.Range("E15").Value = '  -1.95%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '59.134.75'
$ws.Range("E16").Value = '  +0.48%  '

# Row 17
$ws.Range("E17").Value = '  +1.01%  '

# Row 18
$ws.Range("D18").Value = '2.507.65'
$ws.Range("E18").Value = '  -0.43%  '

# Row 19
$ws.Range("D19").Value = '10.98'
$ws.Range("E19").Value = '  -1.88%  '

# Row 20
$ws.Range("E20").Value = '  -0.48%  '

# Row 21
$ws.Range("D21").Value = '322.56'
$ws.Range("E21").Value = '  -0.26%  '

# Row 22
$ws.Range("E22").Value = '  +0.10%  '

# Row 23
$ws.Range("D23").Value = '5.85'
$ws.Range("E23").Value = '  +1.18%  '

# Row 24
$ws.Range("D24").Value = '62.12'
$ws.Range("E24").Value = '  +0.83%  '

# Row 25
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  -2.83%  '

# Row 26
$ws.Range("D26").Value = '0.166'
$ws.Range("E26").Value = '  +1.53%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.66%  '

# Row 28
$ws.Range("D28").Value = '7.82'
$ws.Range("E28").Value = '  +0.55%  '

# Row 29
$ws.Range("D29").Value = '6.92'
$ws.Range("E29").Value = '  +2.71%  '

# Row 30
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0769'
$ws.Range("E30").Value = '  -0.65%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.80'
$ws.Range("E31").Value = '  +0.03%  '

# Row 32
$ws.Range("D32").Value = '161.99'
$ws.Range("E32").Value = '  +2.94%  '

# Row 33
$ws.Range("E33").Value = '  +0.26%  '

# Row 34
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.13'
$ws.Range("E34").Value = '  -6.05%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.45'
$ws.Range("E35").Value = '  +0.56%  '

# Row 36
$ws.Range("D36").Value = '18.45'
$ws.Range("E36").Value = '  -1.06%  '

# Row 37
$ws.Range("D37").Value = '4.25'
$ws.Range("E37").Value = '  -2.51%  '

# Row 38
$ws.Range("E38").Value = '  -1.74%  '

# Row 39
$ws.Range("D39").Value = '36.99'
$ws.Range("E39").Value = '  +0.86%  '

# Row 40
$ws.Range("D40").Value = '3.66'
$ws.Range("E40").Value = '  -0.21%  '

# Row 41
$ws.Range("D41").Value = '0.808'
$ws.Range("E41").Value = '  -3.07%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '5.25'
$ws.Range("E42").Value = '  -8.08%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '281.52'
$ws.Range("E43").Value = '  -6.43%  '

# Row 44
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.54%  '

# Row 45
$ws.Range("D45").Value = '10.85'
$ws.Range("E45").Value = '  +0.59%  '

# Row 46
$ws.Range("E46").Value = '  -0.69%  '

# Row 47
$ws.Range("D47").Value = '0.0932'
$ws.Range("E47").Value = '  +0.15%  '

# Row 48
$ws.Range("D48").Value = '122.47'
$ws.Range("E48").Value = '  -1.76%  '

# Row 49
$ws.Range("E49").Value = '  -0.08%  '

# Row 50
$ws.Range("D50").Value = '0.0514'
$ws.Range("E50").Value = '  -0.38%  '

# Row 51
$ws.Range("D51").Value = '0.0224'
$ws.Range("E51").Value = '  -1.75%  '
